# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row -> new F value for sheet "展览"
$exhibitionUpdates = @{
    7  = 1446
    9  = 115
    14 = 450
    15 = 1389
    17 = 120
    20 = 73
    21 = 662
    24 = 242
    26 = 5976
    28 = 125
    31 = 14674
    32 = 1456
    33 = 228
    36 = 9486
    37 = 645
    38 = 4225
    39 = 161
}

# Map of row -> new F value for sheet "全部类型"
$allTypesUpdates = @{
    7  = 1446
    9  = 115
    14 = 450
    15 = 1389
    17 = 120
    21 = 73
    22 = 662
    26 = 242
    29 = 5976
    31 = 125
    34 = 14674
    35 = 1456
    36 = 228
    39 = 9486
    40 = 645
    41 = 4225
    42 = 161
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
